# Actualizacion automatica del tracker
# Fill in the outcome for row 125 (previously pending) and append three new
# pending rows (126-128) with the latest fetched matches.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    # Force the value to be stored as text (not auto-converted to a date /
    # number by Excel) by using a leading apostrophe, then strip the
    # resulting "quote prefix" style so the cell keeps the default style.
    $ws.Cells.Item($row, $col).Value = "'" + $text
    $ws.Cells.Item($row, $col).Style = "Normal"
}

function Set-EmptyCell($row, $col) {
    # Create an explicit, empty, text-typed cell (mirrors the empty
    # "resultado"/"profit" placeholder cells used for pending matches).
    $ws.Cells.Item($row, $col).Value = "'"
    $ws.Cells.Item($row, $col).Style = "Normal"
}

# --- Row 125: result became known -> fill resultado / profit ---
$ws.Cells.Item(125, 7).Value = "Fallo"
$ws.Cells.Item(125, 8).Value = -1

# --- Row 126: new pending match ---
$ws.Cells.Item(126, 1).Value = 14807100
Set-TextCell 126 2 "2025-10-10"
$ws.Cells.Item(126, 3).Value = "Marco Trungelliti"
$ws.Cells.Item(126, 4).Value = "Carlos Taberner"
$ws.Cells.Item(126, 5).Value = "Gana Carlos Taberner"
$ws.Cells.Item(126, 6).Value = 2.38
Set-EmptyCell 126 7
Set-EmptyCell 126 8

# --- Row 127: new pending match ---
$ws.Cells.Item(127, 1).Value = 14807181
Set-TextCell 127 2 "2025-10-09"
$ws.Cells.Item(127, 3).Value = "Roger Pascual Ferra"
$ws.Cells.Item(127, 4).Value = "Abdullah Shelbayh"
$ws.Cells.Item(127, 5).Value = "Gana Roger Pascual Ferra"
$ws.Cells.Item(127, 6).Value = 3.75
Set-EmptyCell 127 7
Set-EmptyCell 127 8

# --- Row 128: new pending match ---
$ws.Cells.Item(128, 1).Value = 14807179
Set-TextCell 128 2 "2025-10-09"
$ws.Cells.Item(128, 3).Value = "Michael Mmoh"
$ws.Cells.Item(128, 4).Value = "Jay Dylan Hara Friend"
$ws.Cells.Item(128, 5).Value = "Gana Jay Dylan Hara Friend"
$ws.Cells.Item(128, 6).Value = 1.91
Set-EmptyCell 128 7
Set-EmptyCell 128 8
